# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages regeneration @ 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 8314
    $ws.Range("F3").Value = 7742
    $ws.Range("F4").Value = 119
    $ws.Range("F14").Value = 1303
    $ws.Range("F19").Value = 113
}
